$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = 6.2
$ws.Cells.Item(2, 7).Value = 7.4
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(2, 12).Value = 1.01
$ws.Cells.Item(2, 13).Value = 1.01
$ws.Cells.Item(2, 14).Value = 3.85
$ws.Cells.Item(2, 15).Value = 1.3
$ws.Cells.Item(2, 16).Value = 1.98
$ws.Cells.Item(2, 17).Value = 1.81
$ws.Cells.Item(2, 18).Value = 1.38
$ws.Cells.Item(2, 19).Value = 3.05
$ws.Cells.Item(2, 20).Value = 1.75
$ws.Cells.Item(2, 21).Value = 1.78
$ws.Cells.Item(2, 22).Value = 2.48
$ws.Cells.Item(2, 23).Value = 1.15
$ws.Cells.Item(2, 24).Value = 17.5
$ws.Cells.Item(2, 25).Value = 8.800000000000001
$ws.Cells.Item(2, 26).Value = 10
$ws.Cells.Item(2, 27).Value = 16
$ws.Cells.Item(2, 28).Value = 23
$ws.Cells.Item(2, 29).Value = 10.5
$ws.Cells.Item(2, 30).Value = 10.5
$ws.Cells.Item(2, 31).Value = 18
$ws.Cells.Item(2, 32).Value = 60
$ws.Cells.Item(2, 33).Value = 27
$ws.Cells.Item(2, 34).Value = 23
$ws.Cells.Item(2, 35).Value = 38
$ws.Cells.Item(2, 36).Value = 240
$ws.Cells.Item(2, 37).Value = 120
$ws.Cells.Item(2, 38).Value = 110
$ws.Cells.Item(2, 39).Value = 160
$ws.Cells.Item(2, 40).Value = 150
$ws.Cells.Item(2, 41).Value = 9.199999999999999
$ws.Cells.Item(3, 7).Value = 8
$ws.Cells.Item(3, 8).Value = 1.74
$ws.Cells.Item(3, 9).Value = 1.92
$ws.Cells.Item(3, 17).Value = 2.02
$ws.Cells.Item(4, 7).Value = 1.48
$ws.Cells.Item(4, 16).Value = 1.86
$ws.Cells.Item(4, 17).Value = 1.94
$ws.Cells.Item(5, 6).Value = 1.9
$ws.Cells.Item(5, 7).Value = 2.38
$ws.Cells.Item(5, 8).Value = 3.05
$ws.Cells.Item(5, 9).Value = 4.7
$ws.Cells.Item(5, 10).Value = 3.35
$ws.Cells.Item(5, 11).Value = 950
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(6, 7).Value = 8.4
$ws.Cells.Item(6, 8).Value = 1.52
$ws.Cells.Item(6, 9).Value = 1.66
$ws.Cells.Item(6, 11).Value = 5
$ws.Cells.Item(6, 16).Value = 2
$ws.Cells.Item(6, 17).Value = 1.67
$ws.Cells.Item(7, 7).Value = 1.49
$ws.Cells.Item(7, 8).Value = 6.6
$ws.Cells.Item(7, 9).Value = 100
$ws.Cells.Item(7, 10).Value = 4.8
$ws.Cells.Item(7, 11).Value = 15
$ws.Cells.Item(8, 7).Value = 4.7
$ws.Cells.Item(8, 8).Value = 1.81
$ws.Cells.Item(8, 9).Value = 1.84
$ws.Cells.Item(8, 17).Value = 1.78
$ws.Cells.Item(8, 20).Value = 1.75
$ws.Cells.Item(8, 24).Value = 19.5
$ws.Cells.Item(8, 27).Value = 20
$ws.Cells.Item(8, 30).Value = 10
$ws.Cells.Item(8, 36).Value = 120
$ws.Cells.Item(8, 37).Value = 55
$ws.Cells.Item(8, 38).Value = 65
$ws.Cells.Item(8, 39).Value = 110
$ws.Cells.Item(8, 40).Value = 55
$ws.Cells.Item(9, 6).Value = 1.09
$ws.Cells.Item(9, 7).Value = 1.7
$ws.Cells.Item(9, 9).Value = 8.800000000000001
$ws.Cells.Item(10, 7).Value = 3.6
$ws.Cells.Item(10, 10).Value = 3.45
$ws.Cells.Item(10, 17).Value = 2.16
$ws.Cells.Item(10, 25).Value = 9.199999999999999
$ws.Cells.Item(11, 10).Value = 5.9
$ws.Cells.Item(11, 11).Value = 6
$ws.Cells.Item(11, 14).Value = 4.5
$ws.Cells.Item(11, 15).Value = 1.26
$ws.Cells.Item(11, 16).Value = 2.24
$ws.Cells.Item(11, 17).Value = 1.79
$ws.Cells.Item(11, 19).Value = 2.96
$ws.Cells.Item(11, 26).Value = 120
$ws.Cells.Item(11, 28).Value = 8
$ws.Cells.Item(11, 30).Value = 44
$ws.Cells.Item(11, 31).Value = 280
$ws.Cells.Item(11, 38).Value = 44
$ws.Cells.Item(12, 7).Value = 5
$ws.Cells.Item(12, 8).Value = 1.85
$ws.Cells.Item(12, 11).Value = 4
$ws.Cells.Item(12, 16).Value = 1.89
$ws.Cells.Item(12, 20).Value = 1.95
$ws.Cells.Item(12, 21).Value = 1.98
$ws.Cells.Item(12, 24).Value = 13.5
$ws.Cells.Item(12, 25).Value = 8.199999999999999
$ws.Cells.Item(12, 26).Value = 10.5
$ws.Cells.Item(12, 32).Value = 38
$ws.Cells.Item(12, 33).Value = 19.5
$ws.Cells.Item(12, 40).Value = 95
$ws.Cells.Item(12, 41).Value = 13.5
$ws.Cells.Item(13, 6).Value = 2.24
$ws.Cells.Item(13, 7).Value = 2.26
$ws.Cells.Item(13, 8).Value = 3.45
$ws.Cells.Item(13, 9).Value = 3.55
$ws.Cells.Item(13, 12).Value = 1.31
$ws.Cells.Item(13, 13).Value = 1.07
$ws.Cells.Item(13, 22).Value = 1.39
$ws.Cells.Item(13, 23).Value = 1.79
$ws.Cells.Item(13, 32).Value = 15
$ws.Cells.Item(14, 7).Value = 2.86
$ws.Cells.Item(14, 8).Value = 2.64
$ws.Cells.Item(14, 11).Value = 3.9
$ws.Cells.Item(14, 16).Value = 2.08
$ws.Cells.Item(14, 17).Value = 1.67
$ws.Cells.Item(15, 6).Value = 4.8
$ws.Cells.Item(15, 7).Value = 5.5
$ws.Cells.Item(15, 8).Value = 1.72
$ws.Cells.Item(15, 9).Value = 1.82
$ws.Cells.Item(15, 10).Value = 4
$ws.Cells.Item(15, 11).Value = 4.4
$ws.Cells.Item(15, 14).Value = 4.2
$ws.Cells.Item(15, 15).Value = 1.25
$ws.Cells.Item(15, 16).Value = 2.12
$ws.Cells.Item(15, 17).Value = 1.64
$ws.Cells.Item(15, 18).Value = 1.43
$ws.Cells.Item(15, 19).Value = 2.7
$ws.Cells.Item(15, 20).Value = 1.73
$ws.Cells.Item(15, 21).Value = 2.08
$ws.Cells.Item(15, 25).Value = 10.5
$ws.Cells.Item(15, 29).Value = 10
$ws.Cells.Item(15, 36).Value = 140
$ws.Cells.Item(15, 40).Value = 85
$ws.Cells.Item(15, 41).Value = 11.5
$ws.Cells.Item(16, 8).Value = 11
$ws.Cells.Item(16, 10).Value = 5.2
$ws.Cells.Item(16, 14).Value = 4.4
$ws.Cells.Item(16, 15).Value = 1.24
$ws.Cells.Item(16, 16).Value = 2.18
$ws.Cells.Item(16, 18).Value = 1.46
$ws.Cells.Item(16, 19).Value = 2.84
$ws.Cells.Item(16, 20).Value = 2.14
$ws.Cells.Item(16, 23).Value = 3.55
$ws.Cells.Item(16, 24).Value = 24
$ws.Cells.Item(16, 38).Value = 1000
$ws.Cells.Item(17, 8).Value = 28
$ws.Cells.Item(17, 12).Value = 1.15
$ws.Cells.Item(17, 14).Value = 9.199999999999999
$ws.Cells.Item(17, 16).Value = 3.85
$ws.Cells.Item(17, 17).Value = 1.29
$ws.Cells.Item(17, 19).Value = 1.74
$ws.Cells.Item(17, 20).Value = 2.26
$ws.Cells.Item(17, 22).Value = 1.03
$ws.Cells.Item(17, 28).Value = 16
$ws.Cells.Item(17, 29).Value = 1000
$ws.Cells.Item(17, 32).Value = 10.5
$ws.Cells.Item(17, 36).Value = 9.199999999999999
$ws.Cells.Item(17, 40).Value = 2.46
$ws.Cells.Item(18, 6).Value = 2.32
$ws.Cells.Item(18, 7).Value = 2.44
$ws.Cells.Item(18, 8).Value = 3.3
$ws.Cells.Item(18, 14).Value = 3.6
$ws.Cells.Item(18, 15).Value = 1.33
$ws.Cells.Item(18, 16).Value = 1.9
$ws.Cells.Item(18, 17).Value = 1.97
$ws.Cells.Item(18, 19).Value = 3.5
$ws.Cells.Item(18, 20).Value = 1.75
$ws.Cells.Item(18, 21).Value = 2.14
$ws.Cells.Item(18, 23).Value = 1.7
$ws.Cells.Item(18, 24).Value = 16.5
$ws.Cells.Item(18, 28).Value = 10.5
$ws.Cells.Item(18, 29).Value = 8
$ws.Cells.Item(18, 41).Value = 40
$ws.Cells.Item(19, 8).Value = 9.199999999999999
$ws.Cells.Item(19, 14).Value = 5.5
$ws.Cells.Item(19, 16).Value = 2.54
$ws.Cells.Item(19, 20).Value = 1.93
$ws.Cells.Item(19, 32).Value = 8.6
$ws.Cells.Item(19, 36).Value = 11.5
$ws.Cells.Item(20, 9).Value = 1.83
$ws.Cells.Item(20, 12).Value = 1.32
$ws.Cells.Item(20, 16).Value = 2.1
$ws.Cells.Item(20, 21).Value = 2.12
$ws.Cells.Item(20, 22).Value = 2.2
$ws.Cells.Item(20, 25).Value = 9.199999999999999
$ws.Cells.Item(20, 27).Value = 18.5
$ws.Cells.Item(20, 29).Value = 9.4
$ws.Cells.Item(20, 36).Value = 160
$ws.Cells.Item(21, 6).Value = 2
$ws.Cells.Item(21, 8).Value = 3.95
$ws.Cells.Item(21, 9).Value = 4.5
$ws.Cells.Item(21, 11).Value = 3.75
$ws.Cells.Item(21, 16).Value = 1.76
$ws.Cells.Item(21, 17).Value = 2.08
$ws.Cells.Item(22, 16).Value = 1.87
$ws.Cells.Item(23, 6).Value = 4.9
$ws.Cells.Item(23, 8).Value = 1.77
$ws.Cells.Item(23, 9).Value = 1.88
$ws.Cells.Item(23, 16).Value = 1.87
$ws.Cells.Item(24, 6).Value = 1.21
$ws.Cells.Item(24, 9).Value = 18
$ws.Cells.Item(24, 11).Value = 8.4
$ws.Cells.Item(25, 8).Value = 2.92
$ws.Cells.Item(25, 9).Value = 3.15
$ws.Cells.Item(25, 11).Value = 3.45
$ws.Cells.Item(25, 16).Value = 1.68
